# Updates weekly fruit/vegetable price records (Higo - Mercado Mayorista Lo Valledor)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44657
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 13000
$ws.Range("S2").Value = 1857

$ws.Range("D3").Value = 44643
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("S3").Value = 2143

$ws.Range("D4").Value = 44690
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("S4").Value = 1714

$ws.Range("D5").Value = 44312
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range("S5").Value = 1857

$ws.Range("D6").Value = 44312
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 11000
$ws.Range("P6").Value = 11000
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 1571

$ws.Range("D7").Value = 44344
$ws.Range("L7").Value = "Segunda"
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("S7").Value = 1714

$ws.Range("D8").Value = 44307
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 70
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("S8").Value = 2000

$ws.Range("D9").Value = 44307
$ws.Range("L9").Value = "Segunda"
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("S9").Value = 1429

$ws.Range("D10").Value = 44342
$ws.Range("M10").Value = 50

$ws.Range("D11").Value = 44685
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("S11").Value = 2143

$ws.Range("D12").Value = 44685
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("S12").Value = 1714

$ws.Range("D13").Value = 44335
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 14000
$ws.Range("S13").Value = 2000

$ws.Range("D14").Value = 44641
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 13000
$ws.Range("O14").Value = 13000
$ws.Range("P14").Value = 13000
$ws.Range("S14").Value = 1857

$ws.Range("D15").Value = 44644
$ws.Range("M15").Value = 85

$ws.Range("D16").Value = 44694
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("S16").Value = 2143

$ws.Range("D17").Value = 44694
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 75
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("S17").Value = 1714

$ws.Range("D18").Value = 44322
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 11000
$ws.Range("O18").Value = 11000
$ws.Range("P18").Value = 11000
$ws.Range("R18").Value = "Región Metropolitana"
$ws.Range("S18").Value = 1571

$ws.Range("D19").Value = 44306
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 12000
$ws.Range("P19").Value = 12000
$ws.Range("S19").Value = 1714

$ws.Range("D20").Value = 44306
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = 9000
$ws.Range("O20").Value = 9000
$ws.Range("P20").Value = 9000
$ws.Range("S20").Value = 1286

$ws.Range("D21").Value = 44316
$ws.Range("M21").Value = 40
$ws.Range("N21").Value = 13000
$ws.Range("O21").Value = 13000
$ws.Range("P21").Value = 13000
$ws.Range("S21").Value = 1857

$ws.Range("D22").Value = 44316
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = 11000
$ws.Range("O22").Value = 11000
$ws.Range("P22").Value = 11000
$ws.Range("S22").Value = 1571

$ws.Range("D23").Value = 44687
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 15000
$ws.Range("O23").Value = 15000
$ws.Range("P23").Value = 15000
$ws.Range("S23").Value = 2143

$ws.Range("D24").Value = 44687
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 75
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 12000
$ws.Range("S24").Value = 1714

$ws.Range("D25").Value = 44349

$ws.Range("D26").Value = 44302
$ws.Range("M26").Value = 340
$ws.Range("N26").Value = 12000
$ws.Range("P26").Value = 12500
$ws.Range("R26").Value = "Provincia de Santiago"
$ws.Range("S26").Value = 1786

$ws.Range("D27").Value = 44315
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 14000
$ws.Range("O27").Value = 14000
$ws.Range("P27").Value = 14000
$ws.Range("S27").Value = 2000

$ws.Range("D28").Value = 44315
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 80
$ws.Range("O28").Value = 13000
$ws.Range("P28").Value = 12500
$ws.Range("S28").Value = 1786

$ws.Range("D29").Value = 44315
$ws.Range("L29").Value = "Segunda"
$ws.Range("N29").Value = 10000
$ws.Range("O29").Value = 11000
$ws.Range("P29").Value = 10500
$ws.Range("S29").Value = 1500

$ws.Range("D30").Value = 44679
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 150

$ws.Range("D31").Value = 44623
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 30
$ws.Range("N31").Value = 16000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 16000
$ws.Range("S31").Value = 2286

$ws.Range("D32").Value = 44664
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 80
$ws.Range("N32").Value = 14000
$ws.Range("O32").Value = 14000
$ws.Range("P32").Value = 14000
$ws.Range("S32").Value = 2000

$ws.Range("D33").Value = 44664

$ws.Range("D34").Value = 44314
$ws.Range("M34").Value = 20

$ws.Range("D35").Value = 44314
$ws.Range("L35").Value = "Segunda"
$ws.Range("M35").Value = 45
$ws.Range("N35").Value = 11000
$ws.Range("O35").Value = 11000
$ws.Range("P35").Value = 11000
$ws.Range("S35").Value = 1571

$ws.Range("D36").Value = 44300
$ws.Range("M36").Value = 150
$ws.Range("N36").Value = 12000
$ws.Range("O36").Value = 13000
$ws.Range("P36").Value = 12500
$ws.Range("R36").Value = "Provincia de Santiago"
$ws.Range("S36").Value = 1786

$ws.Range("D37").Value = 44321
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 140
$ws.Range("N37").Value = 11000
$ws.Range("P37").Value = 11500
$ws.Range("S37").Value = 1643

$ws.Range("D38").Value = 44321
$ws.Range("M38").Value = 80
$ws.Range("N38").Value = 8000
$ws.Range("O38").Value = 8000
$ws.Range("P38").Value = 8000
$ws.Range("S38").Value = 1143

$ws.Range("D39").Value = 44699
$ws.Range("L39").Value = "Segunda"
$ws.Range("N39").Value = 12000
$ws.Range("O39").Value = 12000
$ws.Range("P39").Value = 12000
$ws.Range("S39").Value = 1714

$ws.Range("D40").Value = 44659
$ws.Range("M40").Value = 50
$ws.Range("N40").Value = 15000
$ws.Range("O40").Value = 15000
$ws.Range("P40").Value = 15000
$ws.Range("S40").Value = 2143

$ws.Range("D41").Value = 44659
$ws.Range("M41").Value = 20
$ws.Range("N41").Value = 12000
$ws.Range("O41").Value = 12000
$ws.Range("P41").Value = 12000
$ws.Range("S41").Value = 1714
